# Regenerate sval data to filter save games: update computed stats for the
# two rows in Sheet1 with new values (TB, d2S, K, IP, sum columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1554434735375247
$ws.Range("C2").Value = 86.29678392075563
$ws.Range("D2").Value = 3.082599426703578
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("G2").Value = 336.5200775150984

# Row 3
$ws.Range("B3").Value = 0.7287194209349384
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 3.082599426703578
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 5.964442013611383
